$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, reusing the existing header style (copy format from H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in new data columns I (I0) and J (IF) for rows 2-9
$data = @(
    @(1, 5),
    @(8, 9),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(6, 8),
    @(4, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
